$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell O2: "Week 4" (matches style of L2/M2/N2 as closely as possible) ---
$h = $ws.Range("O2")
$h.Value = "Week 4"
$h.Font.Name = "Arial"
$h.Font.Size = 12
$h.Font.Bold = $true
$h.Font.Color = 0
$h.Interior.Pattern = 1
$h.Interior.Color = 8630772
$h.HorizontalAlignment = -4108
$h.VerticalAlignment = -4108
$h.WrapText = $true
$h.Borders.Item(7).LineStyle = 1
$h.Borders.Item(7).Weight = 2
$h.Borders.Item(8).LineStyle = 1
$h.Borders.Item(8).Weight = -4138
$h.Borders.Item(9).LineStyle = 1
$h.Borders.Item(9).Weight = 2
$h.Borders.Item(10).LineStyle = 1
$h.Borders.Item(10).Weight = 2

# --- Week 4 notes for Teams 1, 2, 3, 5, 7 (rows 3, 4, 5, 7, 9) ---
$ws.Range("O3").Value = "The team has tested over-sampling and under-sampling and they are getting a very high F1 score and accuracy. So, it's possible that the model is over-fitting. They are investigating to see what could be the problem here."
$ws.Range("O4").Value = "The team is working on using cosine similarity and Jaccard similarity coefficient to find politicians with similar voting patterns."
$ws.Range("O5").Value = "The team has tried various models and getting very low scores ROC-AUC scores. So, they are performing more feature engineering to get better results from the dataset."
$ws.Range("O7").Value = "The team has used cosine similarity on TF-IDF vectors to find similar movies for their recommendation system. The next step will be collaborative filtering in their project. The team is brainstorming ways to validate the goodness of fit for these recommendations."
$ws.Range("O9").Value = "The team has been working through outliers in their dataset. They are now working on their ML models."

# --- New column O width (target raw width 32.41 chars; engine rounds to nearest achievable) ---
$ws.Columns.Item(15).ColumnWidth = 32.41

# --- Update view: scroll so column K is the left-most visible column, row 6 the top, and select O10 ---
$ws.Application.Goto($ws.Range("K6"), $true)
$ws.Range("O10").Select()

Write-Host "done"
